$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# style: horizontal center only (no wrap) - for H22:H24
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").Copy()
$ws.Range("H22:H24").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# style: horizontal center + wrap - for H5:H8, H13:H16, H18:H21
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").WrapText = $true
$ws.Range("Z1").Copy()
$ws.Range("H5:H8").PasteSpecial(-4122)
$ws.Range("H13:H16").PasteSpecial(-4122)
$ws.Range("H18:H21").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# Row 17 - change fill to black across B17:G17
$ws.Range("B17:G17").Interior.Color = 0

Write-Output "Done"
